$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark new instructions as tested ("y") in the relevant checklist columns.
$ws.Range("D4").Value = "y"
$ws.Range("D25").Value = "y"
$ws.Range("D26").Value = "y"
$ws.Range("B29").Value = "y"
$ws.Range("D29").Value = "y"
$ws.Range("B30").Value = "y"

# Update the active selection to D4.
$ws.Range("D4").Select()
